$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.426.08'
$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.869.72'
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.00'
$ws.Range("E5").Value = '  +1.83%  '

$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4739'
$ws.Range("E7").Value = '  +0.34%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2918'
$ws.Range("E8").Value = '  +1.93%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06493'
$ws.Range("E9").Value = '  +0.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.17'
$ws.Range("E10").Value = '  +6.95%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07719'
$ws.Range("E11").Value = '  +0.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.60'
$ws.Range("E12").Value = '  +2.94%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7398'
$ws.Range("E13").Value = '  +5.74%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.873.64'
$ws.Range("E14").Value = '  +0.46%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.142'
$ws.Range("E15").Value = '  +1.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '273.02'
$ws.Range("E16").Value = '  +1.86%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.404.16'
$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.42'
$ws.Range("E18").Value = '  +0.89%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007542'
$ws.Range("E19").Value = '  +0.28%  '

$ws.Range("E20").Value = '  -0.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.120.19'
$ws.Range("E21").Value = '  +0.52%  '

$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.226'
$ws.Range("E23").Value = '  +0.91%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.172'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.286'
$ws.Range("E25").Value = '  -0.18%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.61'
$ws.Range("E26").Value = '  -1.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.79'
$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.923'
$ws.Range("E28").Value = '  +1.01%  '

$ws.Range("E29").Value = '  +1.91%  '

$ws.Range("E30").Value = '  -1.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.505'
$ws.Range("E31").Value = '  -0.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.278'
$ws.Range("E32").Value = '  +1.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.122'
$ws.Range("E33").Value = '  +3.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04827'
$ws.Range("E34").Value = '  +2.42%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.120'
$ws.Range("E35").Value = '  +0.54%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6951'
$ws.Range("E36").Value = '  +1.08%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.718'
$ws.Range("E37").Value = '  +0.49%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01850'
$ws.Range("E38").Value = '  +0.57%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.739'
$ws.Range("E39").Value = '  +0.61%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.311'
$ws.Range("E40").Value = '  -0.12%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.76'
$ws.Range("E41").Value = '  +3.65%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.967'
$ws.Range("E42").Value = '  +4.33%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4189'
$ws.Range("E43").Value = '  +3.31%  '

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.03%  '

$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8378'
$ws.Range("E45").Value = '  -0.15%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.93'
$ws.Range("E46").Value = '  -0.01%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.259'
$ws.Range("E47").Value = '  +0.70%  '

$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.49'
$ws.Range("E48").Value = '  +2.60%  '

$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.985'
$ws.Range("E49").Value = '  -1.15%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '920.24'
$ws.Range("E50").Value = '  -1.53%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05641'
$ws.Range("E51").Value = '  +1.55%  '
